# Rename the first three sheets to carry an "input_" prefix
# (the 4th sheet, component_names, keeps its name).
$wb = $excel.ActiveWorkbook

$wsStoich = $wb.Worksheets.Item(1)   # stoich_coefficients
$wsConc   = $wb.Worksheets.Item(2)   # concentrations
$wsK      = $wb.Worksheets.Item(3)   # k_constants_log10
$wsComp   = $wb.Worksheets.Item(4)   # component_names

$wsStoich.Name = "input_stoich_coefficients"
$wsConc.Name   = "input_concentrations"
$wsK.Name      = "input_k_constants_log10"

# Restore the saved cursor/selection position on each sheet.
$wsStoich.Range("J34").Select() | Out-Null
$wsConc.Range("I33").Select() | Out-Null

# Make the k_constants sheet the active tab with its own selection.
$wsK.Select()
$wsK.Range("L34").Select() | Out-Null
